$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the "last row" (bottom border) formatting from the old last table row (21)
# onto row 17, which will become the new last table row once rows 18-21 are removed.
$ws.Range("B21:J21").Copy() | Out-Null
$ws.Range("B17:J17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Update the "Salario Basico" (G column) values for the remaining worker's two periods.
$ws.Range("G16").Value = 1423500
$ws.Range("G17").Value = 1423500

# Remove the second worker's block (rows 18-21) entirely.
$ws.Range("18:21").Delete() | Out-Null

# Update the summary fields at the top of the sheet.
$ws.Range("E11").Value = 92800
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 2
